$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential notice date from 2021-05-04 to 2021-05-05
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."
$ws.Rows.Item(18).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-15
$ws.Range("D2").Value = 0.05716839648408403
$ws.Range("E2").Value = -0.0001853138753763472
$ws.Range("D3").Value = 0.02341560311863734
$ws.Range("E3").Value = 0.01187648456056989
$ws.Range("D4").Value = 0.03117246859015742
$ws.Range("E4").Value = 0.003823360734085179
$ws.Range("D5").Value = 0.03187896244061512
$ws.Range("E5").Value = 0.02959747434885562
$ws.Range("D6").Value = 0.03712171624373044
$ws.Range("E6").Value = 0.02878852739726034
$ws.Range("D7").Value = 0.0192024631275761
$ws.Range("E7").Value = 0.007965242577842258
$ws.Range("D8").Value = 0.004364105874757908
$ws.Range("E8").Value = -0.001365498406918353
$ws.Range("D9").Value = 0.00696959163066329
$ws.Range("E9").Value = -0.001330039901197066
$ws.Range("D10").Value = 0.07019913591895516
$ws.Range("E10").Value = 0.004527447651386574
$ws.Range("D11").Value = 0.07031831951134727
$ws.Range("E11").Value = 0.003389830508474745
$ws.Range("D12").Value = 0.1477452781778153
$ws.Range("E12").Value = 0.001649218413882103
$ws.Range("D13").Value = 0.3858409892238169
$ws.Range("E13").Value = 0.001225275687029681
$ws.Range("D14").Value = 0.1146029696578438
$ws.Range("E14").Value = 0.009949041494782929
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0.004949429739616873

$ws.Protect()
